$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 63747
$ws.Range("B2").Value = "Maria Liz da Cruz"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Doenca"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45094
$ws.Range("G2").Value = 8755.860000000001

# Row 3
$ws.Range("A3").Value = 19394
$ws.Range("B3").Value = "José Pastor"
$ws.Range("C3").Value = "Financeiro"
$ws.Range("D3").Value = "Viagem de negocios"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45085
$ws.Range("G3").Value = 5384.47

# Row 4
$ws.Range("A4").Value = 82140
$ws.Range("B4").Value = "Pietra Leão"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Doenca"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45097
$ws.Range("G4").Value = 2320.22

# Row 5
$ws.Range("A5").Value = 58592
$ws.Range("B5").Value = "Rael Cavalcante"
$ws.Range("C5").Value = "TI"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 45096
$ws.Range("G5").Value = 2851.41

# Row 6
$ws.Range("A6").Value = 91812
$ws.Range("B6").Value = "Sr. Bryan Mendes"
$ws.Range("C6").Value = "TI"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 3776.62

# Row 7
$ws.Range("A7").Value = 17994
$ws.Range("B7").Value = "Luiz Gustavo Brito"
$ws.Range("C7").Value = "Juridico"
$ws.Range("D7").Value = "Viagem de negocios"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 45079
$ws.Range("G7").Value = 7738.09

# Row 8
$ws.Range("A8").Value = 70745
$ws.Range("B8").Value = "Luiz Fernando da Cruz"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 45093
$ws.Range("G8").Value = 4953.62

# Row 9
$ws.Range("A9").Value = 35678
$ws.Range("B9").Value = "Anthony Gabriel Cardoso"
$ws.Range("C9").Value = "P&D"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45088
$ws.Range("G9").Value = 5695.01

# Row 10
$ws.Range("A10").Value = 49626
$ws.Range("B10").Value = "Milena Melo"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Consulta medica"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45094
$ws.Range("G10").Value = 7780.18

# Row 11
$ws.Range("A11").Value = 3233
$ws.Range("B11").Value = "Maria Flor Guerra"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45103
$ws.Range("G11").Value = 8448.85
